$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D16").Value = 0.25
$ws.Range("C17").Value = 1.5
$ws.Range("D17").Value = 1
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("C20").Value = 2.5
$ws.Range("D20").Value = 0
$ws.Range("G20").Value = 0

$ws.Range("E24").Formula = "=D24-C24"
$ws.Range("E25").Formula = "=D25-C25"
$ws.Range("E26").Formula = "=D26-C26"
